$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 17.01.2022 14:15"

# Row 10 (EuroOil Opuštěná) refreshed values:
#  - new current price in B10 (was old current price)
#  - old price shifted into C10
#  - delta now stored as a literal text string in D10
#  - last-changed timestamp now stored as a literal text string in E10
$ws.Range("B10").Value = 36.4
$ws.Range("C10").Value = 36.2

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "+0.2"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value2 = "2022-01-17 14:15:27"
$ws.Range("E10").Style = "Normal"
